$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colG = @(
    0.2753413333333334,
    0.2753413333333334,
    0.2753413333333334,
    0.2753413333333334,
    0.2753413333333334,
    1.641258,
    1.641258,
    1.641258,
    1.641258,
    1.641258,
    23.78768866666667,
    23.78768866666667,
    23.78768866666667,
    23.78768866666667,
    23.78768866666667,
    0.5982033333333333,
    0.5982033333333333,
    0.5982033333333333,
    0.5982033333333333,
    0.5982033333333333,
    15.22474833333333,
    15.22474833333333,
    15.22474833333333,
    15.22474833333333,
    15.22474833333333
)
for ($i = 0; $i -lt $colG.Length; $i++) {
    $ws.Cells.Item($i + 2, "G").Value2 = $colG[$i]
}

$colH = @(
    0.8260240000000001,
    0.8260240000000001,
    0.8260240000000001,
    0.8260240000000001,
    0.8260240000000001,
    4.923774,
    4.923774,
    4.923774,
    4.923774,
    4.923774,
    71.363066,
    71.363066,
    71.363066,
    71.363066,
    71.363066,
    1.79461,
    1.79461,
    1.79461,
    1.79461,
    1.79461,
    45.674245,
    45.674245,
    45.674245,
    45.674245,
    45.674245
)
for ($i = 0; $i -lt $colH.Length; $i++) {
    $ws.Cells.Item($i + 2, "H").Value2 = $colH[$i]
}

$colI = @(
    0.006630378892106956,
    0.006630378892106956,
    0.006630378892106956,
    0.006630378892106956,
    0.006630378892106956,
    0.03952244389885164,
    0.03952244389885164,
    0.03952244389885164,
    0.03952244389885164,
    0.03952244389885164,
    0.5728213302306416,
    0.5728213302306416,
    0.5728213302306416,
    0.5728213302306416,
    0.5728213302306416,
    0.01440508298011203,
    0.01440508298011203,
    0.01440508298011203,
    0.01440508298011203,
    0.01440508298011203,
    0.3666207639982877,
    0.3666207639982877,
    0.3666207639982877,
    0.3666207639982877,
    0.3666207639982877
)
for ($i = 0; $i -lt $colI.Length; $i++) {
    $ws.Cells.Item($i + 2, "I").Value2 = $colI[$i]
}

$colJ = @(
    0.006630378892106955,
    0.006630378892106955,
    0.006630378892106955,
    0.006630378892106955,
    0.006630378892106955,
    0.03952244389885164,
    0.03952244389885164,
    0.03952244389885164,
    0.03952244389885164,
    0.03952244389885164,
    0.5728213302306416,
    0.5728213302306416,
    0.5728213302306416,
    0.5728213302306416,
    0.5728213302306416,
    0.01440508298011203,
    0.01440508298011203,
    0.01440508298011203,
    0.01440508298011203,
    0.01440508298011203,
    0.3666207639982877,
    0.3666207639982877,
    0.3666207639982877,
    0.3666207639982877,
    0.3666207639982877
)
for ($i = 0; $i -lt $colJ.Length; $i++) {
    $ws.Cells.Item($i + 2, "J").Value2 = $colJ[$i]
}

$colM = @(
    5.375839,
    0.1628146666666667,
    9.994147,
    0.7761303333333333,
    6.522593333333333,
    5.375839,
    0.1628146666666667,
    9.994147,
    0.7761303333333333,
    6.522593333333333,
    5.375839,
    0.1628146666666667,
    9.994147,
    0.7761303333333333,
    6.522593333333333,
    5.375839,
    0.1628146666666667,
    9.994147,
    0.7761303333333333,
    6.522593333333333,
    5.375839,
    0.1628146666666667,
    9.994147,
    0.7761303333333333,
    6.522593333333333
)
for ($i = 0; $i -lt $colM.Length; $i++) {
    $ws.Cells.Item($i + 2, "M").Value2 = $colM[$i]
}

$colN = @(
    16.127517,
    0.488444,
    29.982441,
    2.328391,
    19.56778,
    16.127517,
    0.488444,
    29.982441,
    2.328391,
    19.56778,
    16.127517,
    0.488444,
    29.982441,
    2.328391,
    19.56778,
    16.127517,
    0.488444,
    29.982441,
    2.328391,
    19.56778,
    16.127517,
    0.488444,
    29.982441,
    2.328391,
    19.56778
)
for ($i = 0; $i -lt $colN.Length; $i++) {
    $ws.Cells.Item($i + 2, "N").Value2 = $colN[$i]
}

$colO = @(
    0.2354568587499626,
    0.007131134316291014,
    0.4377345486919088,
    0.03399380269149206,
    0.2856836555503455,
    0.2354568587499626,
    0.007131134316291014,
    0.4377345486919088,
    0.03399380269149206,
    0.2856836555503455,
    0.2354568587499626,
    0.007131134316291014,
    0.4377345486919088,
    0.03399380269149206,
    0.2856836555503455,
    0.2354568587499626,
    0.007131134316291014,
    0.4377345486919088,
    0.03399380269149206,
    0.2856836555503455,
    0.2354568587499626,
    0.007131134316291014,
    0.4377345486919088,
    0.03399380269149206,
    0.2856836555503455
)
for ($i = 0; $i -lt $colO.Length; $i++) {
    $ws.Cells.Item($i + 2, "O").Value2 = $colO[$i]
}

$colP = @(
    0.2354568587499626,
    0.007131134316291014,
    0.4377345486919088,
    0.03399380269149207,
    0.2856836555503455,
    0.2354568587499626,
    0.007131134316291014,
    0.4377345486919088,
    0.03399380269149207,
    0.2856836555503455,
    0.2354568587499626,
    0.007131134316291014,
    0.4377345486919088,
    0.03399380269149207,
    0.2856836555503455,
    0.2354568587499626,
    0.007131134316291014,
    0.4377345486919088,
    0.03399380269149207,
    0.2856836555503455,
    0.2354568587499626,
    0.007131134316291014,
    0.4377345486919088,
    0.03399380269149207,
    0.2856836555503455
)
for ($i = 0; $i -lt $colP.Length; $i++) {
    $ws.Cells.Item($i + 2, "P").Value2 = $colP[$i]
}

$colQ = @(
    1.480190678045334,
    0.04482960740622223,
    2.751801760509334,
    0.2137007608204445,
    1.795939545191111,
    8.823138765462,
    0.267220874184,
    16.402973716926,
    1.273830118626,
    10.70525848908,
    127.8787844541247,
    3.872984601033778,
    237.7376573249007,
    18.46234673408955,
    155.1574195126089,
    3.215844809263333,
    0.09739627631555554,
    5.978532049223333,
    0.464283752501111,
    3.901837073977777,
    81.84579585551833,
    2.478812324975555,
    152.1583728813383,
    11.81638899886611,
    99.30484198067776
)
for ($i = 0; $i -lt $colQ.Length; $i++) {
    $ws.Cells.Item($i + 2, "Q").Value2 = $colQ[$i]
}

$colR = @(
    13.321716102408,
    0.4034664666560001,
    24.766215844584,
    1.923306847384,
    16.16345590672,
    79.408248889158,
    2.404987867656,
    147.626763452334,
    11.464471067634,
    96.34732640172,
    1150.909060087122,
    34.856861409304,
    2139.638915924106,
    166.161120606806,
    1396.41677561348,
    28.94260328337,
    0.87656648684,
    53.80678844301,
    4.17855377251,
    35.1165336658,
    736.6121626996651,
    22.30931092478,
    1369.425355932045,
    106.347500989795,
    893.7435778260999
)
for ($i = 0; $i -lt $colR.Length; $i++) {
    $ws.Cells.Item($i + 2, "R").Value2 = $colR[$i]
}

$colS = @(
    0.001561168186257561,
    0.00004728212244751551,
    0.002902345911992797,
    0.0002253917918281176,
    0.001894190879580965,
    0.00930583049054523,
    0.0002818398559507874,
    0.01730033914326511,
    0.001343518159783127,
    0.01129091624930739,
    0.1348747110410818,
    0.004084865845111195,
    0.2507436864696088,
    0.01947237527733845,
    0.1636456915975014,
    0.003391775588529729,
    0.0001027245815684965,
    0.006305602497168837,
    0.0004896835485804989,
    0.004115296764264471,
    0.0863233734435482,
    0.002614421911213019,
    0.1604825746698733,
    0.01246283391396187,
    0.1047375600596914
)
for ($i = 0; $i -lt $colS.Length; $i++) {
    $ws.Cells.Item($i + 2, "S").Value2 = $colS[$i]
}

$colT = @(
    0.001561168186257561,
    0.0000472821224475155,
    0.002902345911992796,
    0.0002253917918281176,
    0.001894190879580965,
    0.009305830490545232,
    0.0002818398559507874,
    0.01730033914326511,
    0.001343518159783127,
    0.01129091624930739,
    0.1348747110410819,
    0.004084865845111195,
    0.2507436864696088,
    0.01947237527733845,
    0.1636456915975014,
    0.003391775588529729,
    0.0001027245815684965,
    0.006305602497168837,
    0.000489683548580499,
    0.004115296764264471,
    0.0863233734435482,
    0.002614421911213019,
    0.1604825746698733,
    0.01246283391396187,
    0.1047375600596913
)
for ($i = 0; $i -lt $colT.Length; $i++) {
    $ws.Cells.Item($i + 2, "T").Value2 = $colT[$i]
}
